$wb = $excel.ActiveWorkbook

# --- Sheet "Sources" (sheet1) ---
$ws1 = $wb.Worksheets.Item("Sources")
$ws1.Range("B1").Value = "Capex in `$"
$ws1.Range("C1").Value = "Opex in `$/MJ"
$ws1.Range("E1").Value = "CO2 in kg/MJ"

$ws1.Range("A2").Value = "CrOil"
$ws1.Range("B2").Value = 0
$ws1.Range("C2").Value = 0.01
$ws1.Range("D2").Value = "crudeOil"
$ws1.Range("E2").Value = 0.0732

$ws1.Range("A3").Value = "H2"
$ws1.Range("B3").Value = 0
$ws1.Range("C3").Value = 0.05
$ws1.Range("D3").Value = "hydrogen"
$ws1.Range("E3").Value = 0

# --- Sheet "Sinks" (sheet2) ---
$ws2 = $wb.Worksheets.Item("Sinks")
$ws2.Range("B1").Value = "Capex in `$"
$ws2.Range("C1").Value = "Opex in `$/MJ"
$ws2.Range("E1").Value = "Demand in MJ"

$ws2.Range("A2").Value = "Gasoline"
$ws2.Range("B2").Value = 0
$ws2.Range("C2").Value = 0
$ws2.Range("D2").Value = "gasoline"
$ws2.Range("E2").Value = 1000

# --- Sheet "Transformers" (sheet3) ---
$ws3 = $wb.Worksheets.Item("Transformers")
$ws3.Range("C1").Value = "Capex in `$"
$ws3.Range("D1").Value = "Opex in `$/MJ"

$ws3.Range("A2").Value = "Refinery"
$ws3.Range("B2").Value = "crudeoil"
$ws3.Range("C2").Value = 0
$ws3.Range("D2").Value = 0
$ws3.Range("E2").Value = 0.93
$ws3.Range("F2").Value = "gasoline"
$ws3.Range("G2").Value = 1

$ws3.Range("A3").Value = "MtG"
$ws3.Range("B3").Value = "hydrogen"
$ws3.Range("C3").Value = 0
$ws3.Range("D3").Value = 0
$ws3.Range("E3").Value = 0.84
$ws3.Range("F3").Value = "gasoline"
$ws3.Range("G3").Value = 1

# --- Sheet "Connectors" (sheet4) ---
$ws4 = $wb.Worksheets.Item("Connectors")

$ws4.Range("A2").Value = "cr2ref"
$ws4.Range("B2").Value = "CrOil"
$ws4.Range("C2").Value = "Refinery"
$ws4.Range("D2").Value = "crudeoil"

$ws4.Range("A3").Value = "h22mtg"
$ws4.Range("B3").Value = "H2"
$ws4.Range("C3").Value = "MtG"
$ws4.Range("D3").Value = "hydrogen"

$ws4.Range("A4").Value = "ref2gas"
$ws4.Range("B4").Value = "Refinery"
$ws4.Range("C4").Value = "Gasoline"
$ws4.Range("D4").Value = "gasoline"

$ws4.Range("A5").Value = "mtg2gas"
$ws4.Range("B5").Value = "MtG"
$ws4.Range("C5").Value = "Gasoline"
$ws4.Range("D5").Value = "gasoline"

# --- Column widths (bestFit/autofit) to mirror the authored layout ---
$ws1.Range("A1:E3").EntireColumn.AutoFit() | Out-Null
$ws2.Range("A1:E2").EntireColumn.AutoFit() | Out-Null
$ws3.Range("A1:L3").EntireColumn.AutoFit() | Out-Null

# --- Restore selections as recorded in the committed workbook ---
$ws1.Range("C12").Select()
$ws2.Range("E2").Select()
$ws3.Range("D29").Select()
$ws4.Range("E5").Select()
$ws4.Activate()
